$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet becomes the active/selected tab.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate() | Out-Null

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "heading" / "Outstanding" columns one position to the right.
# Excel copies the format of the column to the left (M) onto the newly
# inserted column, including its width, so mirror that explicitly.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $mWidth

# Leave the selection where Excel lands it after the insert.
$ws.Range("R11").Select() | Out-Null
